$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cronograma")

$ws.Range("G16").Value = 1
$ws.Range("G18").Value = 1
$ws.Range("G19").Value = 1
$ws.Range("G20").Value = 1
$ws.Range("G21").Value = 1
$ws.Range("G22").Value = 1
$ws.Range("G23").Value = 1
$ws.Range("G24").Value = 1
$ws.Range("G25").Value = 0.3

$ws.Range("B27").Select()
